# ISAICP-3068: Add 'contact' migration test coverage.
# Adds a new row (row 29) to the "1. Content items" table, duplicating the
# formatting of the last existing row (row 28) and populating it with a new
# "asset_release with contact point" content item.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Expand the worksheet Table (ListObject) by one row; this grows the table
# ref, the AutoFilter ref and the sheet dimension together.
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Copy the formatting (styles, number formats, etc.) of the last populated
# row onto the freshly added row, column by column (column I has no cell on
# row 28, so it is intentionally skipped to avoid materialising it).
$copyCols = @("A","B","C","D","E","F","G","H","J","K","L","M","N","O","P","Q")
foreach ($col in $copyCols) {
    $ws.Range($col + "28").Copy($ws.Range($col + "29"))
}
$ws.Rows.Item(29).RowHeight = $ws.Rows.Item(28).RowHeight()

# Populate the new row's values. Columns are, in order:
# Type of content item | nid | Title of content item | Collection_Name |
# Policy domain 1 | Policy domain 2 | New collection | Migrate | Abstract |
# Logo | Banner | Owner | Collection Owner | Elibrary Creation |
# Pre Moderation | Collection state | Content item state
$ws.Range("A29").Value = "Interoperability Solution"
$ws.Range("B29").Value = 102713
$ws.Range("C29").Value = "asset_release with contact point"
$ws.Range("D29").Value = "Archived collection"
$ws.Range("E29").Value = ""
$ws.Range("F29").Value = "Open government"
$ws.Range("G29").Value = "No"
$ws.Range("H29").Value = "Yes"
$ws.Range("J29").Value = ""
$ws.Range("K29").Value = ""
$ws.Range("L29").Value = ""
$ws.Range("M29").Value = "doe@example.com"
$ws.Range("N29").Value = ""
$ws.Range("O29").Value = ""
$ws.Range("P29").Value = ""
$ws.Range("Q29").Value = ""

# Re-create the "Collection Owner" mailto hyperlink on M29, matching the
# hyperlink already present on M28. Adding a hyperlink re-styles the cell
# with the built-in "Hyperlink" look, so re-apply the original formatting
# afterwards to keep it consistent with the rest of the table.
$ws.Hyperlinks.Add($ws.Range("M29"), "mailto:doe@example.com") | Out-Null
$ws.Range("M28").Copy($ws.Range("M29"))
$ws.Range("M29").Value = "doe@example.com"

# The hidden _FilterDatabase name (created by the sheet's AutoFilter) still
# points at the old range; extend it to cover the new row too.
foreach ($n in $wb.Names) {
    if ($n.Name() -like "*_FilterDatabase*") {
        $n.RefersTo = "='1. Content items'!`$A`$1:`$Q`$29"
    }
}

Write-Host "Added row 29 to Table18915"
